# Adds 24 new rows (218-241) of SXT antibiotic model-summary data to sheet1,
# matching the target diff (dimension grows from A1:H217 to A1:H241).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(218, 1).Value = "(Intercept)"
$ws.Cells.Item(218, 2).Value = [double]"2.461187494463191e+223"
$ws.Cells.Item(218, 3).Value = [double]"165.6636630723051"
$ws.Cells.Item(218, 4).Value = [double]"3.104948364373465"
$ws.Cells.Item(218, 5).Value = [double]"0.001903122447941562"
$ws.Cells.Item(218, 6).Value = [double]"3.51387920061645e+82"
$ws.Cells.Item(218, 7).Value = "Inf"
$ws.Cells.Item(218, 8).Value = "SXT"

$ws.Cells.Item(219, 1).Value = "Year"
$ws.Cells.Item(219, 2).Value = [double]"0.7758372021000897"
$ws.Cells.Item(219, 3).Value = [double]"0.08194274080823544"
$ws.Cells.Item(219, 4).Value = [double]"-3.097438155141306"
$ws.Cells.Item(219, 5).Value = [double]"0.001952011216434471"
$ws.Cells.Item(219, 6).Value = [double]"0.6605317330931308"
$ws.Cells.Item(219, 7).Value = [double]"0.9108309659119398"
$ws.Cells.Item(219, 8).Value = "SXT"

$ws.Cells.Item(220, 1).Value = "Specimen_typeother"
$ws.Cells.Item(220, 2).Value = [double]"0.9383884184111067"
$ws.Cells.Item(220, 3).Value = [double]"0.1829757620131584"
$ws.Cells.Item(220, 4).Value = [double]"-0.3475396023278004"
$ws.Cells.Item(220, 5).Value = [double]"0.7281859720035502"
$ws.Cells.Item(220, 6).Value = [double]"0.6550045061548964"
$ws.Cells.Item(220, 7).Value = [double]"1.342652832176647"
$ws.Cells.Item(220, 8).Value = "SXT"

$ws.Cells.Item(221, 1).Value = "Specimen_typeRespiratory"
$ws.Cells.Item(221, 2).Value = [double]"0.7177329833379086"
$ws.Cells.Item(221, 3).Value = [double]"0.1905516036516234"
$ws.Cells.Item(221, 4).Value = [double]"-1.740513657443875"
$ws.Cells.Item(221, 5).Value = [double]"0.0817688640453153"
$ws.Cells.Item(221, 6).Value = [double]"0.4932518746611949"
$ws.Cells.Item(221, 7).Value = [double]"1.041672522628698"
$ws.Cells.Item(221, 8).Value = "SXT"

$ws.Cells.Item(222, 1).Value = "Specimen_typeUrine"
$ws.Cells.Item(222, 2).Value = [double]"0.7319371569408294"
$ws.Cells.Item(222, 3).Value = [double]"0.1419668590581908"
$ws.Cells.Item(222, 4).Value = [double]"-2.198123012301541"
$ws.Cells.Item(222, 5).Value = [double]"0.02794034103201891"
$ws.Cells.Item(222, 6).Value = [double]"0.5543843945303731"
$ws.Cells.Item(222, 7).Value = [double]"0.967504006125308"
$ws.Cells.Item(222, 8).Value = "SXT"

$ws.Cells.Item(223, 1).Value = "Specimen_typeWound & soft tissues"
$ws.Cells.Item(223, 2).Value = [double]"0.642474533745828"
$ws.Cells.Item(223, 3).Value = [double]"0.1870082084025599"
$ws.Cells.Item(223, 4).Value = [double]"-2.365821815902198"
$ws.Cells.Item(223, 5).Value = [double]"0.01799010046614492"
$ws.Cells.Item(223, 6).Value = [double]"0.4444253285848163"
$ws.Cells.Item(223, 7).Value = [double]"0.9255597801176517"
$ws.Cells.Item(223, 8).Value = "SXT"

$ws.Cells.Item(224, 1).Value = "HospitalCHBH"
$ws.Cells.Item(224, 2).Value = [double]"0.1268359677011731"
$ws.Cells.Item(224, 3).Value = [double]"0.4124061155060416"
$ws.Cells.Item(224, 4).Value = [double]"-5.00686227147491"
$ws.Cells.Item(224, 5).Value = [double]"5.532448383349169e-07"
$ws.Cells.Item(224, 6).Value = [double]"0.05532798351972371"
$ws.Cells.Item(224, 7).Value = [double]"0.2799858889684548"
$ws.Cells.Item(224, 8).Value = "SXT"

$ws.Cells.Item(225, 1).Value = "HospitalCNGMO"
$ws.Cells.Item(225, 2).Value = [double]"1.010065958159261"
$ws.Cells.Item(225, 3).Value = [double]"0.4735704897173304"
$ws.Cells.Item(225, 4).Value = [double]"0.02114919330180446"
$ws.Cells.Item(225, 5).Value = [double]"0.983126643076504"
$ws.Cells.Item(225, 6).Value = [double]"0.3896606230327634"
$ws.Cells.Item(225, 7).Value = [double]"2.557317801074604"
$ws.Cells.Item(225, 8).Value = "SXT"

$ws.Cells.Item(226, 1).Value = "HospitalRabta"
$ws.Cells.Item(226, 2).Value = [double]"0.1764228855845602"
$ws.Cells.Item(226, 3).Value = [double]"0.3519441530864111"
$ws.Cells.Item(226, 4).Value = [double]"-4.929394029525668"
$ws.Cells.Item(226, 5).Value = [double]"8.248506235133874e-07"
$ws.Cells.Item(226, 6).Value = [double]"0.08635695109379382"
$ws.Cells.Item(226, 7).Value = [double]"0.345408436332071"
$ws.Cells.Item(226, 8).Value = "SXT"

$ws.Cells.Item(227, 1).Value = "Ward_ED_ICUED"
$ws.Cells.Item(227, 2).Value = [double]"0.1711157349043254"
$ws.Cells.Item(227, 3).Value = [double]"0.3749884231020588"
$ws.Cells.Item(227, 4).Value = [double]"-4.707919045882528"
$ws.Cells.Item(227, 5).Value = [double]"2.502584441670069e-06"
$ws.Cells.Item(227, 6).Value = [double]"0.080131958378693"
$ws.Cells.Item(227, 7).Value = [double]"0.3505049092046891"
$ws.Cells.Item(227, 8).Value = "SXT"

$ws.Cells.Item(228, 1).Value = "Ward_ED_ICUOther"
$ws.Cells.Item(228, 2).Value = [double]"0.2665860152107615"
$ws.Cells.Item(228, 3).Value = [double]"0.3258540357421799"
$ws.Cells.Item(228, 4).Value = [double]"-4.057210234256819"
$ws.Cells.Item(228, 5).Value = [double]"4.966237475586161e-05"
$ws.Cells.Item(228, 6).Value = [double]"0.1368763152400849"
$ws.Cells.Item(228, 7).Value = [double]"0.4948864600823658"
$ws.Cells.Item(228, 8).Value = "SXT"

$ws.Cells.Item(229, 1).Value = "GenderF"
$ws.Cells.Item(229, 2).Value = [double]"0.6962028365555392"
$ws.Cells.Item(229, 3).Value = [double]"0.09254723088361201"
$ws.Cells.Item(229, 4).Value = [double]"-3.912750557911895"
$ws.Cells.Item(229, 5).Value = [double]"9.125077931006256e-05"
$ws.Cells.Item(229, 6).Value = [double]"0.5804938785028012"
$ws.Cells.Item(229, 7).Value = [double]"0.8344276861652442"
$ws.Cells.Item(229, 8).Value = "SXT"

$ws.Cells.Item(230, 1).Value = "Age_cat0–28 d"
$ws.Cells.Item(230, 2).Value = [double]"1.249911182725254"
$ws.Cells.Item(230, 3).Value = [double]"0.2341104844998627"
$ws.Cells.Item(230, 4).Value = [double]"0.9528513660826736"
$ws.Cells.Item(230, 5).Value = [double]"0.34066538434713"
$ws.Cells.Item(230, 6).Value = [double]"0.7887208235837938"
$ws.Cells.Item(230, 7).Value = [double]"1.976530295857539"
$ws.Cells.Item(230, 8).Value = "SXT"

$ws.Cells.Item(231, 1).Value = "Age_cat29–365 d"
$ws.Cells.Item(231, 2).Value = [double]"1.107987400399107"
$ws.Cells.Item(231, 3).Value = [double]"0.2150207399903162"
$ws.Cells.Item(231, 4).Value = [double]"0.4769084916371542"
$ws.Cells.Item(231, 5).Value = [double]"0.633427284932332"
$ws.Cells.Item(231, 6).Value = [double]"0.7254881178838359"
$ws.Cells.Item(231, 7).Value = [double]"1.686606241698885"
$ws.Cells.Item(231, 8).Value = "SXT"

$ws.Cells.Item(232, 1).Value = "Age_cat1–5 y"
$ws.Cells.Item(232, 2).Value = [double]"1.043253870101134"
$ws.Cells.Item(232, 3).Value = [double]"0.2386728456311455"
$ws.Cells.Item(232, 4).Value = [double]"0.17741670624103"
$ws.Cells.Item(232, 5).Value = [double]"0.8591810857899268"
$ws.Cells.Item(232, 6).Value = [double]"0.6510545340283339"
$ws.Cells.Item(232, 7).Value = [double]"1.661402476121889"
$ws.Cells.Item(232, 8).Value = "SXT"

$ws.Cells.Item(233, 1).Value = "Age_cat6–<30 y"
$ws.Cells.Item(233, 2).Value = [double]"1.241531518511784"
$ws.Cells.Item(233, 3).Value = [double]"0.1788303386410371"
$ws.Cells.Item(233, 4).Value = [double]"1.209781934832335"
$ws.Cells.Item(233, 5).Value = [double]"0.2263625799612052"
$ws.Cells.Item(233, 6).Value = [double]"0.8734142274708139"
$ws.Cells.Item(233, 7).Value = [double]"1.761519814473723"
$ws.Cells.Item(233, 8).Value = "SXT"

$ws.Cells.Item(234, 1).Value = "Age_cat52–<67 y"
$ws.Cells.Item(234, 2).Value = [double]"0.9469174056971396"
$ws.Cells.Item(234, 3).Value = [double]"0.1398685176293984"
$ws.Cells.Item(234, 4).Value = [double]"-0.3899619965735838"
$ws.Cells.Item(234, 5).Value = [double]"0.6965646489739572"
$ws.Cells.Item(234, 6).Value = [double]"0.7198876801532644"
$ws.Cells.Item(234, 7).Value = [double]"1.245945123038501"
$ws.Cells.Item(234, 8).Value = "SXT"

$ws.Cells.Item(235, 1).Value = "Age_cat≥67 y"
$ws.Cells.Item(235, 2).Value = [double]"1.213288963694586"
$ws.Cells.Item(235, 3).Value = [double]"0.1389655615789326"
$ws.Cells.Item(235, 4).Value = [double]"1.391242705974334"
$ws.Cells.Item(235, 5).Value = [double]"0.1641518410995224"
$ws.Cells.Item(235, 6).Value = [double]"0.924503349349912"
$ws.Cells.Item(235, 7).Value = [double]"1.594396457623224"
$ws.Cells.Item(235, 8).Value = "SXT"

$ws.Cells.Item(236, 1).Value = "HospitalCHBH:Ward_ED_ICUED"
$ws.Cells.Item(236, 2).Value = [double]"6.098618702788396"
$ws.Cells.Item(236, 3).Value = [double]"0.5628658481090885"
$ws.Cells.Item(236, 4).Value = [double]"3.212243751234152"
$ws.Cells.Item(236, 5).Value = [double]"0.001317025806803913"
$ws.Cells.Item(236, 6).Value = [double]"2.008071010609528"
$ws.Cells.Item(236, 7).Value = [double]"18.39264582683544"
$ws.Cells.Item(236, 8).Value = "SXT"

$ws.Cells.Item(237, 1).Value = "HospitalCNGMO:Ward_ED_ICUED"
$ws.Cells.Item(237, 8).Value = "SXT"

$ws.Cells.Item(238, 1).Value = "HospitalRabta:Ward_ED_ICUED"
$ws.Cells.Item(238, 2).Value = [double]"3.060845345524598"
$ws.Cells.Item(238, 3).Value = [double]"0.4309699558917862"
$ws.Cells.Item(238, 4).Value = [double]"2.595752022214544"
$ws.Cells.Item(238, 5).Value = [double]"0.009438415636612673"
$ws.Cells.Item(238, 6).Value = [double]"1.331425517903692"
$ws.Cells.Item(238, 7).Value = [double]"7.244033488465726"
$ws.Cells.Item(238, 8).Value = "SXT"

$ws.Cells.Item(239, 1).Value = "HospitalCHBH:Ward_ED_ICUOther"
$ws.Cells.Item(239, 2).Value = [double]"5.469218054929973"
$ws.Cells.Item(239, 3).Value = [double]"0.414385111245283"
$ws.Cells.Item(239, 4).Value = [double]"4.10037814713583"
$ws.Cells.Item(239, 5).Value = [double]"4.124755805203718e-05"
$ws.Cells.Item(239, 6).Value = [double]"2.466092441054754"
$ws.Cells.Item(239, 7).Value = [double]"12.57662936238957"
$ws.Cells.Item(239, 8).Value = "SXT"

$ws.Cells.Item(240, 1).Value = "HospitalCNGMO:Ward_ED_ICUOther"
$ws.Cells.Item(240, 8).Value = "SXT"

$ws.Cells.Item(241, 1).Value = "HospitalRabta:Ward_ED_ICUOther"
$ws.Cells.Item(241, 2).Value = [double]"3.425241855227194"
$ws.Cells.Item(241, 3).Value = [double]"0.3634724295428138"
$ws.Cells.Item(241, 4).Value = [double]"3.387250266642146"
$ws.Cells.Item(241, 5).Value = [double]"0.000705969566382035"
$ws.Cells.Item(241, 6).Value = [double]"1.708151962623899"
$ws.Cells.Item(241, 7).Value = [double]"7.144948620994987"
$ws.Cells.Item(241, 8).Value = "SXT"
